$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 76.18045043945312
$ws.Range("C2").Value = 7.896551609039307
$ws.Range("D2").Value = 40.212406158447266
$ws.Range("H2").Value = 6.25
